$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking text (thousands separators,
# trailing zeros, long decimals) as plain text. Temporarily force those cells to
# Text format while writing so Excel does not re-parse them into floating point
# numbers, then restore the original (General) formatting afterwards.
$priceCellAddrs = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCellAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "31.036.75"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.956.24"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "245.87"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.4903"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").Value = "0.2976"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "0.06870"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "19.15"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "107.90"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").Value = "1.952.19"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "0.07754"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "5.487"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "0.7099"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "283.09"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("D17").Value = "31.063.04"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "13.32"
$ws.Range("D19").Value = "0.000007787"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("D20").Value = "2.207.97"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "5.543"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "6.530"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "9.861"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "169.54"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").Value = "20.07"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "2.238"
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("D29").Value = "0.1058"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").Value = "1.426"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("D31").Value = "1.591"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "4.589"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").Value = "4.501"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").Value = "0.04989"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "0.7612"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("D36").Value = "1.185"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "0.02032"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.569"
$ws.Range("E40").Value = "  +9.97%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "2.170"
$ws.Range("E41").Value = "  +5.85%  "
$ws.Range("D42").Value = "74.50"
$ws.Range("E42").Value = "  +7.53%  "
$ws.Range("D43").Value = "0.4524"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "0.8877"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "8.144"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "980.89"
$ws.Range("E48").Value = "  +8.14%  "
$ws.Range("D49").Value = "0.1272"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "9.421"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "35.87"
$ws.Range("E51").Value = "  +0.61%  "

foreach ($addr in $priceCellAddrs) {
    $ws.Range($addr).ClearFormats()
}
